$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price/Volume columns so numeric-looking
# strings (e.g. "250.68", "1.00", "0.0610") are preserved exactly as text,
# matching the source data which stores these as inline strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "36.498.88"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "1.918.37"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "250.68"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "0.694"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "44.37"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").Value = "58.44"
$ws.Range("E9").Value = "  +9.58%  "
$ws.Range("D10").Value = "0.366"
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "14.44"
$ws.Range("E13").Value = "  +8.28%  "
$ws.Range("D14").Value = "0.801"
$ws.Range("E14").Value = "  +5.88%  "
$ws.Range("D15").Value = "2.197.04"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "5.12"
$ws.Range("E16").Value = "  +4.69%  "
$ws.Range("D17").Value = "1.918.72"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "36.452.14"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").Value = "74.25"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  +4.41%  "
$ws.Range("D21").Value = "250.56"
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("D22").Value = "13.25"
$ws.Range("E22").Value = "  +3.71%  "
$ws.Range("D23").Value = "5.23"
$ws.Range("E23").Value = "  +6.04%  "
$ws.Range("D24").Value = "2.68"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("D27").Value = "167.68"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("D28").Value = "8.79"
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("D29").Value = "18.82"
$ws.Range("E29").Value = "  +2.93%  "
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("D31").Value = "4.54"
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("D32").Value = "0.0610"
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("D33").Value = "1.98"
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("D34").Value = "4.34"
$ws.Range("E34").Value = "  +4.87%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.0881"
$ws.Range("E35").Value = "  +26.56%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -12.73%  "
$ws.Range("D38").Value = "0.864"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("B39").Value = "Gas"
$ws.Range("C39").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D39").Value = "17.13"
$ws.Range("E39").Value = "  +41.08%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.02"
$ws.Range("E40").Value = "  +3.53%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "107.69"
$ws.Range("E41").Value = "  +12.09%  "
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("D43").Value = "17.20"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  +3.46%  "
$ws.Range("D45").Value = "1.348.00"
$ws.Range("E45").Value = "  +3.44%  "
$ws.Range("D46").Value = "2.35"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").Value = "0.0809"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D49").Value = "2.81"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("D51").Value = "2.102.32"
$ws.Range("E51").Value = "  +1.98%  "
